$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: force text type to avoid numeric auto-conversion ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.845.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.374.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.797.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.819.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.374.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0217"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.09"
$ws.Range("D49").Style = "Normal"

# --- Other column updates (B, C, E) — plain text assignments ---
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +11.51%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -1.41%  "
